{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the\n// document's table with its new value (25 cells total). Every \"old\"\n// value below is unique in the document, so a simple ordered\n// search-and-replace over the whole body reproduces the diff exactly,\n// leaving the date header paragraph and all formatting untouched.\nconst replacements = [\n  [\"37\u00d746=\", \"63\u00d742=\"],\n  [\"26\u00d740=\", \"38\u00d747=\"],\n  [\"86\u00d716=\", \"67\u00d757=\"],\n  [\"55\u00d767=\", \"98\u00d727=\"],\n  [\"83\u00d713=\", \"97\u00d798=\"],\n  [\"74\u00d734=\", \"24\u00d754=\"],\n  [\"13\u00d726=\", \"40\u00d740=\"],\n  [\"49\u00d798=\", \"23\u00d753=\"],\n  [\"47\u00d795=\", \"87\u00d772=\"],\n  [\"24\u00d794=\", \"23\u00d786=\"],\n  [\"73\u00d747=\", \"62\u00d734=\"],\n  [\"50\u00d718=\", \"51\u00d715=\"],\n  [\"82\u00d761=\", \"54\u00d748=\"],\n  [\"39\u00d799=\", \"29\u00d784=\"],\n  [\"73\u00d727=\", \"63\u00d755=\"],\n  [\"95\u00d790=\", \"36\u00d789=\"],\n  [\"92\u00d796=\", \"76\u00d773=\"],\n  [\"50\u00d734=\", \"70\u00d747=\"],\n  [\"77\u00d748=\", \"89\u00d732=\"],\n  [\"56\u00d720=\", \"73\u00d771=\"],\n  [\"67\u00d744=\", \"47\u00d783=\"],\n  [\"29\u00d789=\", \"38\u00d775=\"],\n  [\"62\u00d724=\", \"69\u00d714=\"],\n  [\"14\u00d781=\", \"35\u00d759=\"],\n  [\"13\u00d796=\", \"76\u00d758=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt in the\n# document's table with its new value (25 cells total). Every \"old\"\n# value below is unique in the document, so a simple Find/Replace pass\n# over the whole document content reproduces the diff exactly, leaving\n# the date header paragraph and all formatting untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"37\u00d746=\", \"63\u00d742=\"),\n    @(\"26\u00d740=\", \"38\u00d747=\"),\n    @(\"86\u00d716=\", \"67\u00d757=\"),\n    @(\"55\u00d767=\", \"98\u00d727=\"),\n    @(\"83\u00d713=\", \"97\u00d798=\"),\n    @(\"74\u00d734=\", \"24\u00d754=\"),\n    @(\"13\u00d726=\", \"40\u00d740=\"),\n    @(\"49\u00d798=\", \"23\u00d753=\"),\n    @(\"47\u00d795=\", \"87\u00d772=\"),\n    @(\"24\u00d794=\", \"23\u00d786=\"),\n    @(\"73\u00d747=\", \"62\u00d734=\"),\n    @(\"50\u00d718=\", \"51\u00d715=\"),\n    @(\"82\u00d761=\", \"54\u00d748=\"),\n    @(\"39\u00d799=\", \"29\u00d784=\"),\n    @(\"73\u00d727=\", \"63\u00d755=\"),\n    @(\"95\u00d790=\", \"36\u00d789=\"),\n    @(\"92\u00d796=\", \"76\u00d773=\"),\n    @(\"50\u00d734=\", \"70\u00d747=\"),\n    @(\"77\u00d748=\", \"89\u00d732=\"),\n    @(\"56\u00d720=\", \"73\u00d771=\"),\n    @(\"67\u00d744=\", \"47\u00d783=\"),\n    @(\"29\u00d789=\", \"38\u00d775=\"),\n    @(\"62\u00d724=\", \"69\u00d714=\"),\n    @(\"14\u00d781=\", \"35\u00d759=\"),\n    @(\"13\u00d796=\", \"76\u00d758=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
